$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated "SamplesTab" StatQuery (B3): drop the Tumor / Analyte Type columns
# from the SELECT list (all other lines of the query are unchanged).
$newQuery = "SELECT`n" +
            "    DISTINCT (smp.sample_id) AS `"Sample ID`",`n" +
            "    sp.participant_id AS `"Participant ID`", `n" +
            "    s.study_name AS `"Study Name`",`n" +
            "    s.phs_accession AS Accession`n" +
            "FROM `n" +
            "    df_participant sp`n" +
            "JOIN `n" +
            "    df_study s ON sp.`"study.phs_accession`" = s.phs_accession`n" +
            "JOIN `n" +
            "    df_sample smp ON smp.`"participant.study_participant_id`" = sp.study_participant_id`n" +
            "JOIN`n" +
            "    df_diagnosis d ON d.`"participant.study_participant_id`" = sp.study_participant_id`n" +
            "JOIN`n" +
            "    df_program p ON p.program_acronym = s.`"program.program_acronym`"`n" +
            "JOIN`n" +
            "    df_file f1 ON f1.`"sample.sample_id`" = smp.sample_id`n" +
            "JOIN`n" +
            "    df_genomic_info gi ON gi.`"file.file_id`" = f1.file_id`n" +
            "WHERE `n" +
            "  s.phs_accession = 'phs001819' AND gi.library_selection = 'Random'`n" +
            "ORDER BY `n" +
            "    smp.sample_id ASC`n" +
            "LIMIT 100;"

$ws.Range("B3").Value = $newQuery

# Move the view/selection up to the Participants/Samples rows, matching the
# saved cursor position in the edited workbook.
$ws.Range("B3").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
